$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: new credentials pair (vj2805 / vishalj2805)
$ws.Range("A4").Value = "vj2805"
$ws.Range("B4").Value = "vishalj2805"

# Row 9: new credentials pair (vsj / vsj2805)
$ws.Range("A9").Value = "vsj"
$ws.Range("B9").Value = "vsj2805"

# Select B9 as the active cell, matching the final selection state
$ws.Range("B9").Select()
